# The sheet contains only a header row (A1:G1). The edit adds a new data
# row (row 2) with the year value back into column A, and leaves the
# selection on the newly written cell - matching what Excel would do
# after typing a value into A2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2005

$ws.Range("A2").Select()
